# Update the skos:concept references to the correctly-cased skos:Concept
# in column D of Sheet1 (generated from the updated Google Sheet / .ttl export).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$oldValue = "iop:VariableSet,skos:concept"
$newValue = "iop:VariableSet,skos:Concept"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
